$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column K (11) so its stored width goes from 116 to 107
$ws.Columns.Item(11).ColumnWidth = 106.17

# Rows where the instructor name (once spaces are stripped) actually matches
# Sofía's name -> mark the comparison as VERDADERO
$verdaderoRows = @(9, 13, 15, 18, 19, 20, 24)
foreach ($r in $verdaderoRows) {
    $ws.Cells.Item($r, 11).Value = "VERDADERO"
}

# Row 23: the instructor-side query name was missing a trailing "S" in
# "HILAMOS" -> fix the name and refresh the discrepancy message accordingly
$ws.Cells.Item(23, 9).Value = "YULIANDREACABALHILAMOS"
$ws.Cells.Item(23, 11).Value = "FALSO - Discrepancia en Nombre: Instructores (YULIANDREACABALHILAMOS) vs Sofía (YULI ANDREA CABAL HILAMO)"
